# Standardized colors and adding wide data
$wb = $excel.ActiveWorkbook

# 1) Rename "Preventative Health" -> "Prevention" wherever it occurs on Sheet1
$ws1 = $wb.Worksheets.Item("Sheet1")
$used = $ws1.UsedRange
$found = $used.Find("Preventative Health")
if ($found) {
    $firstAddress = $found.Address()
    do {
        $found.Value = "Prevention"
        $found = $used.FindNext($found)
    } while ($found -and ($found.Address() -ne $firstAddress))
}

# 2) Remove the empty Sheet2 and Sheet3
$excel.DisplayAlerts = $false
foreach ($name in @("Sheet2", "Sheet3")) {
    foreach ($sheet in $wb.Worksheets) {
        if ($sheet.Name -eq $name) {
            [void]$sheet.Delete()
            break
        }
    }
}
$excel.DisplayAlerts = $true
